$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.807.93"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "1.708.35"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9922"
$ws.Range("E4").Value = "  -1.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.56"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9951"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3910"
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4072"
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.495"
$ws.Range("E9").Value = "  -1.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "54.47"
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9921"
$ws.Range("E11").Value = "  -1.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08826"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.22"
$ws.Range("E13").Value = "  +10.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.484"
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.156"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001363"
$ws.Range("E16").Value = "  +3.25%  "
$ws.Range("D17").Value = "1.701.24"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07171"
$ws.Range("E19").Value = "  +1.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.68"
$ws.Range("E20").Value = "  +4.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.318"
$ws.Range("E21").Value = "  +3.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9948"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.38"
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("D24").Value = "24.799.08"
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("E25").Value = "  -3.24%  "
$ws.Range("E26").Value = "  -0.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.07"
$ws.Range("E27").Value = "  +1.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.19"
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.933"
$ws.Range("E29").Value = "  +15.32%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.488"
$ws.Range("E30").Value = "  -3.63%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "143.82"
$ws.Range("E31").Value = "  +6.15%  "
$ws.Range("D32").Value = "1.889.23"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08831"
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.175"
$ws.Range("E34").Value = "  +10.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.064"
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.290"
$ws.Range("E36").Value = "  -4.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.03111"
$ws.Range("E37").Value = "  +4.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8792"
$ws.Range("E38").Value = "  +14.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2815"
$ws.Range("E39").Value = "  +2.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.96"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09199"
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.29"
$ws.Range("E42").Value = "  -0.78%  "
$ws.Range("E43").Value = "  +1.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.30"
$ws.Range("E44").Value = "  +8.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7557"
$ws.Range("E45").Value = "  +5.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.686"
$ws.Range("E46").Value = "  +3.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.253"
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.404"
$ws.Range("E48").Value = "  +4.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9943"
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08261"
$ws.Range("E51").Value = "  +3.58%  "
